$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---

# Overview sheet: columns E (zh-cn) and F (de-de), rows 2-4
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

# zh-cn sheet: column C (Status), rows 2-4
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

# de-de sheet: column C (Status), rows 2-4
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- Narrow the Status columns to fit new text ---
# Target raw OOXML column width is 13.4101845877511 (as produced by the
# report-generation tool). Excel's ColumnWidth (character-width) property
# is quantized to 1/6ths of a character when round-tripped through COM, so
# the closest attainable value is used here (ColumnWidth=12.5 -> stored
# width 13.333333333333334, nearest reachable value to 13.4101845877511).
$targetColumnWidth = 12.5
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
